# Update "想去人数" (number of people interested) counts for three events
# that each gained one more interested person, on both the "展览" sheet
# and the mirrored "全部类型" sheet.
#   F3: 1691 -> 1692
#   F4:   27 -> 28
#   F9:  611 -> 612

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1692
    $ws.Range("F4").Value = 28
    $ws.Range("F9").Value = 612
}
